# Weekly price-sheet update: insert a new observation row above the
# current row 168 (shifting the remaining records down by one, which the
# diff shows rippling all the way through row 253 -> 254), then populate
# the freshly inserted row with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 168; rows 168:253 shift down to 169:254.
$ws.Rows.Item(168).Insert()

# Populate the newly inserted row 168 with the new weekly record.
$ws.Range("A168").Value = 6
$ws.Range("B168").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C168").Value = "Metropolitana"
$ws.Range("D168").Value = 44813
$ws.Range("E168").Value = 13
$ws.Range("F168").Value = 100112022
$ws.Range("G168").Value = "Arveja Verde"
$ws.Range("H168").Value = "Perfection"
$ws.Range("I168").Value = "Primera"
$ws.Range("J168").Value = 220
$ws.Range("K168").Value = 35000
$ws.Range("L168").Value = 37000
$ws.Range("M168").Value = 35909
$ws.Range("N168").Value = "$/malla 25 kilos"
$ws.Range("O168").Value = "Provincia de Huasco"
$ws.Range("P168").Value = 1436
$ws.Range("Q168").Value = 25
$ws.Range("R168").Value = "Hortaliza"
